$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "823 HK Equity"
$ws.Range("A4").Value = "6823 HK Equity"

$ws.Range("B3").Value = "LINK REITs"

$ws.Range("D3").Value = "HK"

$ws.Range("B4").Value = "HKT Trust and HKT Ltd"

$ws.Range("D4").Value = "HK"

$ws.Range("E3").Value = "Equity, Listed equities"
$ws.Range("E4").Value = "Equity, Listed equities"

$ws.Range("F2").Select()
